$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "cap nhat luong thang 9" - update the September salary sheet:
# - Thuy (row 9) left / her row is removed, Dung/Nam shift up one row
# - Hai (row 3) now gets a bonus marker referencing Thuy's 0.7 rate
# - Tuan (row 8) now gets a bonus marker referencing Duong's 0.7 rate

$ws.Range("K3").Value = "*|Thuy:0.7"
$ws.Range("K8").Value = "*|Duong:0.7"

$ws.Rows("9").Delete()

# Column K now holds the longer "*|Name:0.7" markers - widen it to fit.
$ws.Columns("11:11").ColumnWidth = 11

$ws.Range("L8").Select() | Out-Null
